$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 held summary statistics for column B (years of data).
# D18 used to compute the average deviation; now it reports the median.
$ws.Range("D18").Formula = "=MEDIAN(B1:B18)"

# E18 used to hold the population standard deviation formula; it now
# just holds a plain value, and the stdev/variance formulas shift one
# column to the right (E->F, F->G).
$ws.Range("E18").Value = 0
$ws.Range("F18").Formula = "=STDEV.P(B1:B18)"
$ws.Range("G18").Formula = "=VAR.P(B1:B18)"

# H18 (VAR.S) is unchanged.

$ws.Range("H18").Select()
